# Updated cryptos list with latest price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'24.626.32"
$ws.Range("E2").Value = "'  -0.33%  "

# Row 3
$ws.Range("D3").Value = "'1.672.80"
$ws.Range("E3").Value = "'  -0.71%  "

# Row 4
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "'  +0.47%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'307.39"
$ws.Range("E5").Value = "'  +0.31%  "

# Row 6
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  +0.56%  "

# Row 7
$ws.Range("D7").Value = "'0.3687"
$ws.Range("E7").Value = "'  -0.07%  "

# Row 8
$ws.Range("D8").Value = "'48.13"
$ws.Range("E8").Value = "'  +0.27%  "

# Row 9
$ws.Range("D9").Value = "'0.3370"
$ws.Range("E9").Value = "'  -1.55%  "

# Row 10
$ws.Range("D10").Value = "'1.176"
$ws.Range("E10").Value = "'  +1.05%  "

# Row 11
$ws.Range("D11").Value = "'0.07325"
$ws.Range("E11").Value = "'  +1.31%  "

# Row 12
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "'  +0.50%  "

# Row 13
$ws.Range("D13").Value = "'6.181"
$ws.Range("E13").Value = "'  +1.37%  "

# Row 14
$ws.Range("D14").Value = "'20.52"
$ws.Range("E14").Value = "'  +1.79%  "

# Row 15
$ws.Range("D15").Value = "'6.790"
$ws.Range("E15").Value = "'  +1.43%  "

# Row 16
$ws.Range("D16").Value = "'1.671.86"
$ws.Range("E16").Value = "'  -0.50%  "

# Row 17
$ws.Range("D17").Value = "'0.00001097"

# Row 18
$ws.Range("D18").Value = "'0.06647"
$ws.Range("E18").Value = "'  -0.09%  "

# Row 19
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "'  +0.58%  "

# Row 20
$ws.Range("D20").Value = "'81.59"
$ws.Range("E20").Value = "'  +0.97%  "

# Row 21
$ws.Range("D21").Value = "'16.81"
$ws.Range("E21").Value = "'  +2.53%  "

# Row 22
$ws.Range("E22").Value = "'  +2.17%  "

# Row 23
$ws.Range("D23").Value = "'12.67"
$ws.Range("E23").Value = "'  +4.68%  "

# Row 24
$ws.Range("D24").Value = "'24.571.69"
$ws.Range("E24").Value = "'  -0.22%  "

# Row 25
$ws.Range("E25").Value = "'  +1.11%  "

# Row 26
$ws.Range("D26").Value = "'2.680"
$ws.Range("E26").Value = "'  +0.55%  "

# Row 27
$ws.Range("D27").Value = "'19.87"
$ws.Range("E27").Value = "'  +2.11%  "

# Row 28
$ws.Range("D28").Value = "'149.04"
$ws.Range("E28").Value = "'  -2.49%  "

# Row 29
$ws.Range("D29").Value = "'130.08"
$ws.Range("E29").Value = "'  +1.77%  "

# Row 30
$ws.Range("D30").Value = "'1.855.33"
$ws.Range("E30").Value = "'  -0.63%  "

# Row 31
$ws.Range("D31").Value = "'1.218"
$ws.Range("E31").Value = "'  +24.74%  "

# Row 32
$ws.Range("D32").Value = "'6.516"
$ws.Range("E32").Value = "'  +4.01%  "

# Row 33
$ws.Range("D33").Value = "'4.150"
$ws.Range("E33").Value = "'  +3.11%  "

# Row 34
$ws.Range("D34").Value = "'0.08599"
$ws.Range("E34").Value = "'  +2.38%  "

# Row 35
$ws.Range("D35").Value = "'13.28"
$ws.Range("E35").Value = "'  +7.51%  "

# Row 36
$ws.Range("D36").Value = "'1.710"
$ws.Range("E36").Value = "'  +1.52%  "

# Row 37
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'5.413"
$ws.Range("E37").Value = "'  +1.99%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06448"
$ws.Range("E38").Value = "'  +1.73%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02350"
$ws.Range("E39").Value = "'  +2.00%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.819"
$ws.Range("E40").Value = "'  +2.24%  "

# Row 41
$ws.Range("D41").Value = "'0.2164"
$ws.Range("E41").Value = "'  +3.88%  "

# Row 42
$ws.Range("D42").Value = "'1.230"
$ws.Range("E42").Value = "'  -0.33%  "

# Row 43
$ws.Range("D43").Value = "'0.6248"
$ws.Range("E43").Value = "'  +2.57%  "

# Row 44
$ws.Range("D44").Value = "'1.002"
$ws.Range("E44").Value = "'  +0.67%  "

# Row 45
$ws.Range("D45").Value = "'13.25"
$ws.Range("E45").Value = "'  +2.15%  "

# Row 46
$ws.Range("D46").Value = "'3.776"
$ws.Range("E46").Value = "'  +0.47%  "

# Row 47
$ws.Range("D47").Value = "'0.5926"
$ws.Range("E47").Value = "'  +1.09%  "

# Row 48
$ws.Range("D48").Value = "'2.044"
$ws.Range("E48").Value = "'  +2.15%  "

# Row 49
$ws.Range("D49").Value = "'126.07"
$ws.Range("E49").Value = "'  +0.55%  "

# Row 50
$ws.Range("D50").Value = "'0.07115"
$ws.Range("E50").Value = "'  -1.43%  "

# Row 51
$ws.Range("D51").Value = "'76.87"
$ws.Range("E51").Value = "'  +1.63%  "
